# Update "想去人数" (column F) figures on both the "展览" and "全部类型"
# sheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new value } for column F updates.
$updates = @{
    "展览" = @{
        3  = 1857
        4  = 486
        7  = 2525
        8  = 164
        9  = 87
        11 = 1517
        12 = 524
        15 = 227
        18 = 206
        20 = 216
        21 = 11
        22 = 163
        23 = 52
        24 = 1608
        26 = 391
        27 = 569
        30 = 410
    }
    "全部类型" = @{
        3  = 1857
        5  = 486
        8  = 2525
        9  = 164
        10 = 87
        12 = 1517
        13 = 524
        16 = 227
        19 = 206
        21 = 216
        22 = 11
        23 = 163
        24 = 52
        25 = 1608
        27 = 391
        28 = 569
        31 = 410
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowsMap[$row]
    }
}

$wb.Save()
